# Generate Report for Handoff
#
# Two files have just been handed off for localization:
#   305cf013-4591-4861-bf7c-846e8d1b4bc3.md
#   545787c5-b7d6-44f5-b9e9-de658f4829c9.md
#
# They need a new row in every sheet of the report (Overview, zh-cn, de-de),
# inserted right after the already-handed-back a8e87e34 entry and before the
# pre-existing edcda5fe "Ready for handoff" entry (which stays, just moves
# down to the bottom of the table).

$wb = $excel.ActiveWorkbook

function Remove-HyperlinkAt($ws, [string]$addr) {
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq $addr) {
            $h.Delete()
        }
    }
}

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")
$loOv = $wsOv.ListObjects.Item(1)

# Stash the current last row (edcda5fe) before we start shuffling rows
$ov_old = @()
for ($c = 1; $c -le 7; $c++) {
    $ov_old += $wsOv.Cells.Item(3, $c).Value()
}

Remove-HyperlinkAt $wsOv '$B$3'

# Row 3 becomes the newly-handed-off 305cf013 entry
$wsOv.Cells.Item(3,1).Value = "305cf013-4591-4861-bf7c-846e8d1b4bc3.md"
$wsOv.Cells.Item(3,2).Value = "e2e\305cf013-4591-4861-bf7c-846e8d1b4bc3.md"
$wsOv.Cells.Item(3,3).Value = ".md"
$wsOv.Cells.Item(3,4).Value = ""
$wsOv.Cells.Item(3,5).Value = "Ready for handoff"
$wsOv.Cells.Item(3,6).Value = "Ready for handoff"
$wsOv.Cells.Item(3,7).Value = "2016-09-01 12:45:30"
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/305cf013459148614861bf7c846e8d1b4bc3000/e2e/305cf013-4591-4861-bf7c-846e8d1b4bc3.md", "", "", "e2e\305cf013-4591-4861-bf7c-846e8d1b4bc3.md") | Out-Null

# Add a row for the second newly-handed-off entry (545787c5) -> row 4
$loOv.ListRows.Add() | Out-Null
$wsOv.Cells.Item(4,1).Value = "545787c5-b7d6-44f5-b9e9-de658f4829c9.md"
$wsOv.Cells.Item(4,2).Value = "e2e\545787c5-b7d6-44f5-b9e9-de658f4829c9.md"
$wsOv.Cells.Item(4,3).Value = ".md"
$wsOv.Cells.Item(4,4).Value = ""
$wsOv.Cells.Item(4,5).Value = "Ready for handoff"
$wsOv.Cells.Item(4,6).Value = "Ready for handoff"
$wsOv.Cells.Item(4,7).Value = "2016-09-01 12:45:30"
$wsOv.Cells.Item(4,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOv.Hyperlinks.Add($wsOv.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/545787c5b7d644f5b9e9de658f4829c900000000/e2e/545787c5-b7d6-44f5-b9e9-de658f4829c9.md", "", "", "e2e\545787c5-b7d6-44f5-b9e9-de658f4829c9.md") | Out-Null

# Add a row to host the previously-last (edcda5fe) entry -> row 5
$loOv.ListRows.Add() | Out-Null
for ($c = 1; $c -le 7; $c++) {
    $wsOv.Cells.Item(5, $c).Value = $ov_old[$c-1]
}
$wsOv.Cells.Item(5,7).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOv.Hyperlinks.Add($wsOv.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44802276168e6755566efba8d54cbb9827fd62ee/e2e/edcda5fe-0528-488d-a511-47b278f2f23e.md", "", "", "e2e\edcda5fe-0528-488d-a511-47b278f2f23e.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)

$zh_old = @()
for ($c = 1; $c -le 16; $c++) {
    $zh_old += $wsZh.Cells.Item(3, $c).Value()
}

Remove-HyperlinkAt $wsZh '$A$3'

$wsZh.Cells.Item(3,1).Value  = "305cf013-4591-4861-bf7c-846e8d1b4bc3.md"
$wsZh.Cells.Item(3,2).Value  = ".md"
$wsZh.Cells.Item(3,3).Value  = "Ready for handoff"
$wsZh.Cells.Item(3,4).Value  = "e2e"
$wsZh.Cells.Item(3,5).Value  = "ht"
$wsZh.Cells.Item(3,6).Value  = "False"
$wsZh.Cells.Item(3,7).Value  = "305cf013-4591-4861-bf7c-846e8d1b4bc3.d89046c0306b30a67654e59368ff9093633169c7.zh-cn.xlf"
$wsZh.Cells.Item(3,8).Value  = "2016-09-01 12:45:25"
$wsZh.Cells.Item(3,9).Value  = ""
$wsZh.Cells.Item(3,10).Value = ""
$wsZh.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(3,12).Value = ""
$wsZh.Cells.Item(3,13).Value = "True"
$wsZh.Cells.Item(3,14).Value = ""
$wsZh.Cells.Item(3,15).Value = "False"
$wsZh.Cells.Item(3,16).Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/305cf013459148614861bf7c846e8d1b4bc3000/e2e/305cf013-4591-4861-bf7c-846e8d1b4bc3.md", "", "", "305cf013-4591-4861-bf7c-846e8d1b4bc3.md") | Out-Null

$loZh.ListRows.Add() | Out-Null
$wsZh.Cells.Item(4,1).Value  = "545787c5-b7d6-44f5-b9e9-de658f4829c9.md"
$wsZh.Cells.Item(4,2).Value  = ".md"
$wsZh.Cells.Item(4,3).Value  = "Ready for handoff"
$wsZh.Cells.Item(4,4).Value  = "e2e"
$wsZh.Cells.Item(4,5).Value  = "ht"
$wsZh.Cells.Item(4,6).Value  = "False"
$wsZh.Cells.Item(4,7).Value  = "545787c5-b7d6-44f5-b9e9-de658f4829c9.89cdad1631365539faba68adbcf2d747b4aaf222.zh-cn.xlf"
$wsZh.Cells.Item(4,8).Value  = "2016-09-01 12:45:25"
$wsZh.Cells.Item(4,9).Value  = ""
$wsZh.Cells.Item(4,10).Value = ""
$wsZh.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(4,12).Value = ""
$wsZh.Cells.Item(4,13).Value = "True"
$wsZh.Cells.Item(4,14).Value = ""
$wsZh.Cells.Item(4,15).Value = "False"
$wsZh.Cells.Item(4,16).Value = ""
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/545787c5b7d644f5b9e9de658f4829c900000000/e2e/545787c5-b7d6-44f5-b9e9-de658f4829c9.md", "", "", "545787c5-b7d6-44f5-b9e9-de658f4829c9.md") | Out-Null

$loZh.ListRows.Add() | Out-Null
for ($c = 1; $c -le 16; $c++) {
    $wsZh.Cells.Item(5, $c).Value = $zh_old[$c-1]
}
$wsZh.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44802276168e6755566efba8d54cbb9827fd62ee/e2e/edcda5fe-0528-488d-a511-47b278f2f23e.md", "", "", "edcda5fe-0528-488d-a511-47b278f2f23e.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)

$de_old = @()
for ($c = 1; $c -le 16; $c++) {
    $de_old += $wsDe.Cells.Item(3, $c).Value()
}

Remove-HyperlinkAt $wsDe '$A$3'

$wsDe.Cells.Item(3,1).Value  = "305cf013-4591-4861-bf7c-846e8d1b4bc3.md"
$wsDe.Cells.Item(3,2).Value  = ".md"
$wsDe.Cells.Item(3,3).Value  = "Ready for handoff"
$wsDe.Cells.Item(3,4).Value  = "e2e"
$wsDe.Cells.Item(3,5).Value  = "ht"
$wsDe.Cells.Item(3,6).Value  = "False"
$wsDe.Cells.Item(3,7).Value  = "305cf013-4591-4861-bf7c-846e8d1b4bc3.d89046c0306b30a67654e59368ff9093633169c7.de-de.xlf"
$wsDe.Cells.Item(3,8).Value  = "2016-09-01 12:45:30"
$wsDe.Cells.Item(3,9).Value  = ""
$wsDe.Cells.Item(3,10).Value = ""
$wsDe.Cells.Item(3,11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(3,12).Value = ""
$wsDe.Cells.Item(3,13).Value = "True"
$wsDe.Cells.Item(3,14).Value = ""
$wsDe.Cells.Item(3,15).Value = "False"
$wsDe.Cells.Item(3,16).Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/305cf013459148614861bf7c846e8d1b4bc3000/e2e/305cf013-4591-4861-bf7c-846e8d1b4bc3.md", "", "", "305cf013-4591-4861-bf7c-846e8d1b4bc3.md") | Out-Null

$loDe.ListRows.Add() | Out-Null
$wsDe.Cells.Item(4,1).Value  = "545787c5-b7d6-44f5-b9e9-de658f4829c9.md"
$wsDe.Cells.Item(4,2).Value  = ".md"
$wsDe.Cells.Item(4,3).Value  = "Ready for handoff"
$wsDe.Cells.Item(4,4).Value  = "e2e"
$wsDe.Cells.Item(4,5).Value  = "ht"
$wsDe.Cells.Item(4,6).Value  = "False"
$wsDe.Cells.Item(4,7).Value  = "545787c5-b7d6-44f5-b9e9-de658f4829c9.89cdad1631365539faba68adbcf2d747b4aaf222.de-de.xlf"
$wsDe.Cells.Item(4,8).Value  = "2016-09-01 12:45:30"
$wsDe.Cells.Item(4,9).Value  = ""
$wsDe.Cells.Item(4,10).Value = ""
$wsDe.Cells.Item(4,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item(4,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(4,12).Value = ""
$wsDe.Cells.Item(4,13).Value = "True"
$wsDe.Cells.Item(4,14).Value = ""
$wsDe.Cells.Item(4,15).Value = "False"
$wsDe.Cells.Item(4,16).Value = ""
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/545787c5b7d644f5b9e9de658f4829c900000000/e2e/545787c5-b7d6-44f5-b9e9-de658f4829c9.md", "", "", "545787c5-b7d6-44f5-b9e9-de658f4829c9.md") | Out-Null

$loDe.ListRows.Add() | Out-Null
for ($c = 1; $c -le 16; $c++) {
    $wsDe.Cells.Item(5, $c).Value = $de_old[$c-1]
}
$wsDe.Cells.Item(5,8).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Cells.Item(5,11).NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/44802276168e6755566efba8d54cbb9827fd62ee/e2e/edcda5fe-0528-488d-a511-47b278f2f23e.md", "", "", "edcda5fe-0528-488d-a511-47b278f2f23e.md") | Out-Null

Write-Host "Handoff rows added for 305cf013 and 545787c5 across Overview/zh-cn/de-de."
